# Remove the rows for years 2007, 2008 and 2009 (original rows 2-4).
# Excel shifts all subsequent rows up, so the former row 5 (2010年)
# becomes row 2, and the former row 13 (2018年) becomes row 10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:4").Delete()
